# Generate Report for Archive
#
# 1) Status text changed from "Ready for handoff" to "In Translation"
#    (appears in the per-language Status columns on every sheet).
# 2) The Status column(s) were narrowed (from a stored width of
#    ~17.216 down to ~13.410).

$wb = $excel.ActiveWorkbook

$newStatus = "In Translation"

$overview = $wb.Worksheets.Item("Overview")
$zhcn     = $wb.Worksheets.Item("zh-cn")
$dede     = $wb.Worksheets.Item("de-de")

# --- Update the status text wherever it appears --------------------------
# Touch the known cells directly (rather than scanning every used cell) so
# the shared-string table for every other, untouched cell ("True"/"False"/
# dates/etc.) is left completely alone.
$overview.Range("E2").Value = $newStatus   # zh-cn status
$overview.Range("F2").Value = $newStatus   # de-de status
$zhcn.Range("C2").Value = $newStatus       # Status column
$dede.Range("C2").Value = $newStatus       # Status column

# --- Narrow the Status column(s) ------------------------------------------
# The stored OOXML column width shrinks from 17.2159881591797 to
# 13.4101845877511. This engine's ColumnWidth setter (like real Excel's)
# only accepts values on a discrete character/pixel grid, so 12.5 is the
# closest settable width that lands on that target.
$newColumnWidth = 12.5

$overview.Columns(5).ColumnWidth = $newColumnWidth  # zh-cn status column
$overview.Columns(6).ColumnWidth = $newColumnWidth  # de-de status column
$zhcn.Columns(3).ColumnWidth = $newColumnWidth       # Status column
$dede.Columns(3).ColumnWidth = $newColumnWidth       # Status column
